# Refresh Sheets market-price columns (H:N) via scheduled runner update.
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5754676.5
$ws.Range("I62").Value = 6417713.5
$ws.Range("J62").Value = 8356.333000000001
$ws.Range("K62").Value = 6417713.5
$ws.Range("L62").Value = 8356.333000000001
$ws.Range("M62").Value = -6417089.5
$ws.Range("N62").Value = -9604.333000000001
$ws.Range("H65").Value = 5754676.5
$ws.Range("I65").Value = 6417713.5
$ws.Range("J65").Value = 8356.333000000001
$ws.Range("K65").Value = 32088567.5
$ws.Range("L65").Value = 41781.665
$ws.Range("M65").Value = -32085447.5
$ws.Range("N65").Value = -48021.665
$ws.Range("H125").Value = 2493.8667
$ws.Range("I125").Value = 3090.4443
$ws.Range("J125").Value = 1599
$ws.Range("K125").Value = 27813.9987
$ws.Range("L125").Value = 14391
$ws.Range("M125").Value = -25353.9987
$ws.Range("N125").Value = -19311
$ws.Range("H129").Value = 37038450
$ws.Range("I129").Value = 1339.4
$ws.Range("J129").Value = 83334830
$ws.Range("K129").Value = 4018.2
$ws.Range("L129").Value = 250004490
$ws.Range("M129").Value = 981.7999999999997
$ws.Range("N129").Value = -250014490
$ws.Range("H138").Value = 4466.26
$ws.Range("I138").Value = 3686.077
$ws.Range("J138").Value = 4624.7344
$ws.Range("K138").Value = 11058.231
$ws.Range("L138").Value = 13874.2032
$ws.Range("M138").Value = -5918.231
$ws.Range("N138").Value = -24154.2032

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 2532500
$ws.Range("I34").Value = 2532500
$ws.Range("K34").Value = 2532500
$ws.Range("M34").Value = -2532229
$ws.Range("H40").Value = 31000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20411034
$ws.Range("I20").Value = 35716308
$ws.Range("K20").Value = 35716308
$ws.Range("M20").Value = -35716061
$ws.Range("H94").Value = 5563.49
$ws.Range("I94").Value = 4776.472
$ws.Range("K94").Value = 4776.472
$ws.Range("M94").Value = -4325.472
$ws.Range("H137").Value = 59000
$ws.Range("J137").Value = 59000
$ws.Range("L137").Value = 59000
$ws.Range("N137").Value = -69200

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8635.272000000001
$ws.Range("I31").Value = 1998.75
$ws.Range("J31").Value = 9550.655000000001
$ws.Range("K31").Value = 1998.75
$ws.Range("L31").Value = 9550.655000000001
$ws.Range("M31").Value = -1703.75
$ws.Range("N31").Value = -10140.655
$ws.Range("H34").Value = 8635.272000000001
$ws.Range("I34").Value = 1998.75
$ws.Range("J34").Value = 9550.655000000001
$ws.Range("K34").Value = 1998.75
$ws.Range("L34").Value = 9550.655000000001
$ws.Range("M34").Value = -1796.75
$ws.Range("N34").Value = -9954.655000000001
$ws.Range("H58").Value = 6592.5
$ws.Range("I58").Value = 5284.579
$ws.Range("J58").Value = 8249.200000000001
$ws.Range("K58").Value = 5284.579
$ws.Range("L58").Value = 8249.200000000001
$ws.Range("M58").Value = -5081.579
$ws.Range("N58").Value = -8655.200000000001
$ws.Range("H99").Value = 3476374
$ws.Range("I99").Value = 5559266
$ws.Range("J99").Value = 4887.1665
$ws.Range("K99").Value = 5559266
$ws.Range("L99").Value = 4887.1665
$ws.Range("M99").Value = -5557768
$ws.Range("N99").Value = -7883.1665
$ws.Range("H122").Value = 5798
$ws.Range("I122").Value = 6469.5
$ws.Range("J122").Value = 5462.25
$ws.Range("K122").Value = 19408.5
$ws.Range("L122").Value = 16386.75
$ws.Range("M122").Value = -16958.5
$ws.Range("N122").Value = -21286.75
$ws.Range("H126").Value = 3476374
$ws.Range("I126").Value = 5559266
$ws.Range("J126").Value = 4887.1665
$ws.Range("K126").Value = 16677798
$ws.Range("L126").Value = 14661.4995
$ws.Range("M126").Value = -16675328
$ws.Range("N126").Value = -19601.4995
$ws.Range("H132").Value = 9228.35
$ws.Range("I132").Value = 8834.9
$ws.Range("K132").Value = 26504.7
$ws.Range("M132").Value = -23974.7
$ws.Range("H136").Value = 6592.5
$ws.Range("I136").Value = 5284.579
$ws.Range("J136").Value = 8249.200000000001
$ws.Range("K136").Value = 15853.737
$ws.Range("L136").Value = 24747.6
$ws.Range("M136").Value = -13303.737
$ws.Range("N136").Value = -29847.6
$ws.Range("H141").Value = 202857.31
$ws.Range("J141").Value = 206530.78
$ws.Range("L141").Value = 206530.78
$ws.Range("N141").Value = -216890.78

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 8951
$ws.Range("I81").Value = 2999.5
$ws.Range("J81").Value = 12918.667
$ws.Range("K81").Value = 8998.5
$ws.Range("L81").Value = 38756.001
$ws.Range("M81").Value = -7875.5
$ws.Range("N81").Value = -41002.001
$ws.Range("H84").Value = 8951
$ws.Range("I84").Value = 2999.5
$ws.Range("J84").Value = 12918.667
$ws.Range("K84").Value = 26995.5
$ws.Range("L84").Value = 116268.003
$ws.Range("M84").Value = -21379.5
$ws.Range("N84").Value = -127500.003

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6265.8667
$ws.Range("I70").Value = 5855.25
$ws.Range("J70").Value = 6539.6113
$ws.Range("K70").Value = 5855.25
$ws.Range("L70").Value = 6539.6113
$ws.Range("M70").Value = -5585.25
$ws.Range("N70").Value = -7079.6113
$ws.Range("H73").Value = 6265.8667
$ws.Range("I73").Value = 5855.25
$ws.Range("J73").Value = 6539.6113
$ws.Range("K73").Value = 5855.25
$ws.Range("L73").Value = 6539.6113
$ws.Range("M73").Value = -4919.25
$ws.Range("N73").Value = -8411.6113
$ws.Range("H80").Value = 71439384
$ws.Range("I80").Value = 200005840
$ws.Range("J80").Value = 13572.223
$ws.Range("K80").Value = 200005840
$ws.Range("L80").Value = 13572.223
$ws.Range("M80").Value = -200004842
$ws.Range("N80").Value = -15568.223
$ws.Range("H83").Value = 71439384
$ws.Range("I83").Value = 200005840
$ws.Range("J83").Value = 13572.223
$ws.Range("K83").Value = 1000029200
$ws.Range("L83").Value = 67861.11500000001
$ws.Range("M83").Value = -1000024208
$ws.Range("N83").Value = -77845.11500000001
$ws.Range("H97").Value = 1466.6842
$ws.Range("I97").Value = 1504.6
$ws.Range("K97").Value = 1504.6
$ws.Range("M97").Value = -1008.6
$ws.Range("H122").Value = 6169.448
$ws.Range("I122").Value = 4977.4287
$ws.Range("K122").Value = 14932.2861
$ws.Range("M122").Value = -12482.2861

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 42270.5
$ws.Range("J50").Value = 44542
$ws.Range("L50").Value = 44542
$ws.Range("N50").Value = -45816
$ws.Range("H68").Value = 1655.4445
$ws.Range("J68").Value = 2100
$ws.Range("L68").Value = 2100
$ws.Range("N68").Value = -3598
$ws.Range("H71").Value = 1655.4445
$ws.Range("J71").Value = 2100
$ws.Range("L71").Value = 10500
$ws.Range("N71").Value = -17988
$ws.Range("H82").Value = 5266.5713
$ws.Range("I82").Value = 1379.25
$ws.Range("K82").Value = 1379.25
$ws.Range("M82").Value = -1018.25
$ws.Range("H85").Value = 5266.5713
$ws.Range("I85").Value = 1379.25
$ws.Range("K85").Value = 1379.25
$ws.Range("M85").Value = -131.25
$ws.Range("H136").Value = 7479.7437
$ws.Range("I136").Value = 8230.714
$ws.Range("J136").Value = 7059.2
$ws.Range("K136").Value = 24692.142
$ws.Range("L136").Value = 21177.6
$ws.Range("M136").Value = -22142.142
$ws.Range("N136").Value = -26277.6

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 42999.668
$ws.Range("J42").Value = 42999.668
$ws.Range("L42").Value = 42999.668
$ws.Range("N42").Value = -43755.668
$ws.Range("H49").Value = 49299.5
$ws.Range("I49").Value = 46759.4
$ws.Range("K49").Value = 46759.4
$ws.Range("M49").Value = -46529.4
$ws.Range("H100").Value = 1824.5862
$ws.Range("I100").Value = 1276.4348
$ws.Range("K100").Value = 2552.8696
$ws.Range("M100").Value = -2011.8696
$ws.Range("H125").Value = 70999.664
$ws.Range("J125").Value = 70999.664
$ws.Range("L125").Value = 70999.664
$ws.Range("N125").Value = -80839.664

# ARM row 40: LeveProfitNQ (M40) cleared - no longer applicable
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M40").ClearContents()
